$wb = $excel.ActiveWorkbook

# Helper: force a literal text value into a cell (even when it looks like
# a number, e.g. "004685" or "5.00") without leaving behind a stray
# number-format style. We build a tiny text formula, then flatten it to a
# plain value in place.
function Set-TextValue($rng, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)  # xlPasteValues
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right before the "总计" sheet. Clone
#    "2021-Q2" (same per-quarter layout/formatting) and overwrite its
#    content, rather than building a blank sheet from scratch, so the
#    sheet keeps the same sheetPr/format boilerplate as its siblings.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2021-Q2")
$totalSheet = $wb.Worksheets.Item("总计")

$q2.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("2021-Q2 (2)")
$newSheet.Name = "2022-Q1"

# Header labels (D1 changes from "基金金额" to "基金规模"; the rest repeat).
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Index cell + data row.
$newSheet.Range("A2").Value = 0

Set-TextValue $newSheet.Range("B2") "004685"
Set-TextValue $newSheet.Range("C2") "金元顺安元启灵活配置混合"
Set-TextValue $newSheet.Range("D2") "5.00"
Set-TextValue $newSheet.Range("E2") "75.79"
Set-TextValue $newSheet.Range("F2") "0.98"
Set-TextValue $newSheet.Range("G2") "0.0490"

$newSheet.Range("H2").Value = 10

# ---------------------------------------------------------------------
# 2. Prepend a new row to the "总计" sheet for 2022-Q1, shifting the
#    existing rows down (done manually, cell-by-cell, so no extra
#    formatting/styles get dragged in along the way).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# old row 3 (2021-Q1) -> row 4
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$totalSheet.Range("A4").Value = 2
Set-TextValue $totalSheet.Range("B4") "2021-Q1"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.03

# old row 2 (2021-Q2) -> row 3
$totalSheet.Range("A2:D2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$totalSheet.Range("A3").Value = 1
Set-TextValue $totalSheet.Range("B3") "2021-Q2"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.03

# new row 2 (2022-Q1)
$totalSheet.Range("A2").Value = 0
Set-TextValue $totalSheet.Range("B2") "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.05

Write-Output "done"
